$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Fall")
$ws.Range("B2").Value = 72
$ws.Range("B3").Value = 74
$ws.Range("B4").Value = 77
$ws.Range("B5").Value = 80
$ws.Range("B6").Value = 83
$ws.Range("B14").Value = 106
$ws.Range("B15").Value = 109
$ws.Range("B16").Value = 112
$ws.Range("B21").Value = 126
$ws.Range("B22").Value = 129

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Range("B3").Value = 67
$ws.Range("B4").Value = 70
$ws.Range("B5").Value = 73
$ws.Range("B6").Value = 76
$ws.Range("B10").Value = 88
$ws.Range("B13").Value = 96
$ws.Range("B14").Value = 99
$ws.Range("B15").Value = 102
$ws.Range("B16").Value = 105
$ws.Range("B17").Value = 108
$ws.Range("B18").Value = 110
$ws.Range("B19").Value = 113
$ws.Range("B20").Value = 116
$ws.Range("B21").Value = 119
$ws.Range("B22").Value = 122

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Range("B2").Value = 56
$ws.Range("B3").Value = 59
$ws.Range("B4").Value = 61
$ws.Range("B5").Value = 64
$ws.Range("B6").Value = 67
$ws.Range("B7").Value = 70
$ws.Range("B8").Value = 73
$ws.Range("B9").Value = 76
$ws.Range("B10").Value = 79
$ws.Range("B11").Value = 82
$ws.Range("B12").Value = 84
$ws.Range("B13").Value = 87
$ws.Range("B14").Value = 90
$ws.Range("B15").Value = 93
$ws.Range("B19").Value = 104
$ws.Range("B25").Value = 122
$ws.Range("B26").Value = 125
$ws.Range("B27").Value = 127

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Range("B2").Value = 72
$ws.Range("B3").Value = 74
$ws.Range("B4").Value = 77
$ws.Range("B5").Value = 80
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 86
$ws.Range("B8").Value = 89
$ws.Range("B9").Value = 92
$ws.Range("B10").Value = 94
$ws.Range("B11").Value = 97
$ws.Range("B12").Value = 100
$ws.Range("B13").Value = 103
$ws.Range("B14").Value = 106
$ws.Range("B15").Value = 109
$ws.Range("B16").Value = 112
$ws.Range("B17").Value = 114
$ws.Range("B18").Value = 117
$ws.Range("B19").Value = 120
$ws.Range("B20").Value = 123
$ws.Range("B21").Value = 126
$ws.Range("B22").Value = 129
$ws.Range("B23").Value = 130
$ws.Range("B24").Value = 130
$ws.Range("B25").Value = 130
$ws.Range("B26").Value = 130
$ws.Range("B27").Value = 130
$ws.Range("B28").Value = 130
$ws.Range("B29").Value = 130

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Range("B2").Value = 44
$ws.Range("B3").Value = 47
$ws.Range("B4").Value = 49
$ws.Range("B5").Value = 52
$ws.Range("B6").Value = 55
$ws.Range("B7").Value = 58
$ws.Range("B8").Value = 61
$ws.Range("B9").Value = 64
$ws.Range("B10").Value = 67
$ws.Range("B11").Value = 70
$ws.Range("B13").Value = 75
$ws.Range("B14").Value = 78
$ws.Range("B15").Value = 81
$ws.Range("B16").Value = 84
$ws.Range("B17").Value = 87
$ws.Range("B18").Value = 90
$ws.Range("B19").Value = 92
$ws.Range("B20").Value = 95
$ws.Range("B21").Value = 98
$ws.Range("B22").Value = 101
$ws.Range("B23").Value = 104
$ws.Range("B24").Value = 107
$ws.Range("B25").Value = 110
$ws.Range("B26").Value = 112
$ws.Range("B27").Value = 115
$ws.Range("B28").Value = 118
$ws.Range("B29").Value = 121
$ws.Range("B30").Value = 124

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Range("B4").Value = 41
$ws.Range("B5").Value = 44
$ws.Range("B6").Value = 47
$ws.Range("B7").Value = 50
$ws.Range("B8").Value = 53
$ws.Range("B9").Value = 56
$ws.Range("B10").Value = 59
$ws.Range("B11").Value = 62
$ws.Range("B19").Value = 84
$ws.Range("B30").Value = 116
$ws.Range("B31").Value = 119
$ws.Range("B32").Value = 122
$ws.Range("B33").Value = 125
